$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds opportunity IDs as text (matches source data which is not numeric-typed)
$ws.Range("A2:A6").NumberFormat = "@"

$ws.Range("A2").Value = "1330202"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1330202"
$ws.Range("C2").Value = "Multimedia Ai Designer"
$ws.Range("D2").Value = "Belgrade, Serbia"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "3 applicants"
$ws.Range("G2").Value = "9 - 12 Weeks"
$ws.Range("H2").Value = "Beyond Creative"

$ws.Range("A3").Value = "1328300"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328300"
$ws.Range("C3").Value = "content creator"
$ws.Range("D3").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "15 applicants"
$ws.Range("G3").Value = "9 - 12 Weeks"
$ws.Range("H3").Value = "Markit"

$ws.Range("A4").Value = "1326963"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326963"
$ws.Range("C4").Value = "Deutsch Ambassador and Instructor"
$ws.Range("D4").Value = "Sousse, Tunisia"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "1 applicant"
$ws.Range("G4").Value = "9 - 12 Weeks"
$ws.Range("H4").Value = "Progress Professional Center"

$ws.Range("A5").Value = "1311536"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1311536"
$ws.Range("C5").Value = "Accelerate Romania | Managing Co-founder"
$ws.Range("D5").Value = "Bucharest, Romania"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "156 applicants"
$ws.Range("G5").Value = "9 - 12 Weeks"
$ws.Range("H5").Value = "Skulptor"

$ws.Range("A6").Value = "1301518"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1301518"
$ws.Range("C6").Value = "MARKETING"
$ws.Range("D6").Value = "Yıldırım, Türkiye"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "75 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "OMTEC Automotive"

# ColumnWidth (character units) maps to a slightly larger stored "width" in the
# XML (stored = ColumnWidth + 11/12, then rounded to the nearest 1/12 character
# by Excel's pixel-based quantization). Back the requested value off by that
# fixed offset (plus a small safety margin) so the saved <col width="..."/>
# lands exactly on the target integer values from the diff.
$ws.Columns.Item(3).ColumnWidth = 43 - 11/12 + 0.02
$ws.Columns.Item(4).ColumnWidth = 70 - 11/12 + 0.02
$ws.Columns.Item(6).ColumnWidth = 17 - 11/12 + 0.02
$ws.Columns.Item(8).ColumnWidth = 31 - 11/12 + 0.02
